$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new blank rows before the current row 325, pushing the existing
# rows 325-329 down to become rows 331-335.
$ws.Range("A325:A330").EntireRow.Insert()

# ---------------------------------------------------------------------
# New row 325 (was row 325 data, now overwritten with updated values)
# ---------------------------------------------------------------------
$ws.Cells.Item(325, 1).Value2 = 3
$ws.Cells.Item(325, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(325, 3).Value2 = "Coquimbo"
$ws.Cells.Item(325, 4).Value2 = 44578
$ws.Cells.Item(325, 5).Value2 = 5
$ws.Cells.Item(325, 6).Value2 = 100112027
$ws.Cells.Item(325, 7).Value2 = "Melón"
$ws.Cells.Item(325, 8).Value2 = "Calameño"
$ws.Cells.Item(325, 9).Value2 = "Extra"
$ws.Cells.Item(325, 10).Value2 = 1100
$ws.Cells.Item(325, 11).Value2 = 1200
$ws.Cells.Item(325, 12).Value2 = 1200
$ws.Cells.Item(325, 13).Value2 = 1200
$ws.Cells.Item(325, 14).Value2 = "$/unidad"
$ws.Cells.Item(325, 15).Value2 = "Provincia de Talca"
$ws.Cells.Item(325, 16).Value2 = 1200
$ws.Cells.Item(325, 17).Value2 = 1
$ws.Cells.Item(325, 18).Value2 = "Hortaliza"

# ---------------------------------------------------------------------
# New row 326
# ---------------------------------------------------------------------
$ws.Cells.Item(326, 1).Value2 = 3
$ws.Cells.Item(326, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(326, 3).Value2 = "Coquimbo"
$ws.Cells.Item(326, 4).Value2 = 44578
$ws.Cells.Item(326, 5).Value2 = 5
$ws.Cells.Item(326, 6).Value2 = 100112027
$ws.Cells.Item(326, 7).Value2 = "Melón"
$ws.Cells.Item(326, 8).Value2 = "Calameño"
$ws.Cells.Item(326, 9).Value2 = "Primera"
$ws.Cells.Item(326, 10).Value2 = 1200
$ws.Cells.Item(326, 11).Value2 = 800
$ws.Cells.Item(326, 12).Value2 = 800
$ws.Cells.Item(326, 13).Value2 = 800
$ws.Cells.Item(326, 14).Value2 = "$/unidad"
$ws.Cells.Item(326, 15).Value2 = "Provincia de Talca"
$ws.Cells.Item(326, 16).Value2 = 800
$ws.Cells.Item(326, 17).Value2 = 1
$ws.Cells.Item(326, 18).Value2 = "Hortaliza"

# ---------------------------------------------------------------------
# New row 327
# ---------------------------------------------------------------------
$ws.Cells.Item(327, 1).Value2 = 3
$ws.Cells.Item(327, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(327, 3).Value2 = "Coquimbo"
$ws.Cells.Item(327, 4).Value2 = 44578
$ws.Cells.Item(327, 5).Value2 = 5
$ws.Cells.Item(327, 6).Value2 = 100112027
$ws.Cells.Item(327, 7).Value2 = "Melón"
$ws.Cells.Item(327, 8).Value2 = "Calameño"
$ws.Cells.Item(327, 9).Value2 = "Segunda"
$ws.Cells.Item(327, 10).Value2 = 1100
$ws.Cells.Item(327, 11).Value2 = 500
$ws.Cells.Item(327, 12).Value2 = 500
$ws.Cells.Item(327, 13).Value2 = 500
$ws.Cells.Item(327, 14).Value2 = "$/unidad"
$ws.Cells.Item(327, 15).Value2 = "Provincia de Talca"
$ws.Cells.Item(327, 16).Value2 = 500
$ws.Cells.Item(327, 17).Value2 = 1
$ws.Cells.Item(327, 18).Value2 = "Hortaliza"

# ---------------------------------------------------------------------
# New row 328
# ---------------------------------------------------------------------
$ws.Cells.Item(328, 1).Value2 = 3
$ws.Cells.Item(328, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(328, 3).Value2 = "Coquimbo"
$ws.Cells.Item(328, 4).Value2 = 44578
$ws.Cells.Item(328, 5).Value2 = 5
$ws.Cells.Item(328, 6).Value2 = 100112027
$ws.Cells.Item(328, 7).Value2 = "Melón"
$ws.Cells.Item(328, 8).Value2 = "Tuna"
$ws.Cells.Item(328, 9).Value2 = "Extra"
$ws.Cells.Item(328, 10).Value2 = 950
$ws.Cells.Item(328, 11).Value2 = 1200
$ws.Cells.Item(328, 12).Value2 = 1200
$ws.Cells.Item(328, 13).Value2 = 1200
$ws.Cells.Item(328, 14).Value2 = "$/unidad"
$ws.Cells.Item(328, 15).Value2 = "Provincia de Talca"
$ws.Cells.Item(328, 16).Value2 = 1200
$ws.Cells.Item(328, 17).Value2 = 1
$ws.Cells.Item(328, 18).Value2 = "Hortaliza"

# ---------------------------------------------------------------------
# New row 329 (brand-new record)
# ---------------------------------------------------------------------
$ws.Cells.Item(329, 1).Value2 = 3
$ws.Cells.Item(329, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(329, 3).Value2 = "Coquimbo"
$ws.Cells.Item(329, 4).Value2 = 44578
$ws.Cells.Item(329, 5).Value2 = 5
$ws.Cells.Item(329, 6).Value2 = 100112027
$ws.Cells.Item(329, 7).Value2 = "Melón"
$ws.Cells.Item(329, 8).Value2 = "Tuna"
$ws.Cells.Item(329, 9).Value2 = "Primera"
$ws.Cells.Item(329, 10).Value2 = 900
$ws.Cells.Item(329, 11).Value2 = 800
$ws.Cells.Item(329, 12).Value2 = 800
$ws.Cells.Item(329, 13).Value2 = 800
$ws.Cells.Item(329, 14).Value2 = "$/unidad"
$ws.Cells.Item(329, 15).Value2 = "Provincia de Talca"
$ws.Cells.Item(329, 16).Value2 = 800
$ws.Cells.Item(329, 17).Value2 = 1
$ws.Cells.Item(329, 18).Value2 = "Hortaliza"

# ---------------------------------------------------------------------
# New row 330 (brand-new record)
# ---------------------------------------------------------------------
$ws.Cells.Item(330, 1).Value2 = 3
$ws.Cells.Item(330, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(330, 3).Value2 = "Coquimbo"
$ws.Cells.Item(330, 4).Value2 = 44578
$ws.Cells.Item(330, 5).Value2 = 5
$ws.Cells.Item(330, 6).Value2 = 100112027
$ws.Cells.Item(330, 7).Value2 = "Melón"
$ws.Cells.Item(330, 8).Value2 = "Tuna"
$ws.Cells.Item(330, 9).Value2 = "Segunda"
$ws.Cells.Item(330, 10).Value2 = 850
$ws.Cells.Item(330, 11).Value2 = 500
$ws.Cells.Item(330, 12).Value2 = 500
$ws.Cells.Item(330, 13).Value2 = 500
$ws.Cells.Item(330, 14).Value2 = "$/unidad"
$ws.Cells.Item(330, 15).Value2 = "Provincia de Talca"
$ws.Cells.Item(330, 16).Value2 = 500
$ws.Cells.Item(330, 17).Value2 = 1
$ws.Cells.Item(330, 18).Value2 = "Hortaliza"

Write-Output ("Dimension now: " + $ws.UsedRange.Rows.Count)
